$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.826168894767761
$ws.Range("B1").Value = 4.6402907371521
$ws.Range("C1").Value = 4.000645160675049
$ws.Range("D1").Value = 1.230011701583862
$ws.Range("E1").Value = 0.7893243432044983
